$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled data
$ws.Range("F2").Value = -6
$ws.Range("F6").Value = -5
$ws.Range("F7").Value = 0
$ws.Range("F11").Value = 3
$ws.Range("F12").Value = -1
$ws.Range("F13").Value = 2
$ws.Range("F14").Value = 5
$ws.Range("F20").Value = -4
